$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.115.11"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.051.57"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.33"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.33"
$ws.Range("E8").Value = "  -3.45%  "
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.109"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.94"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.874"
$ws.Range("E13").Value = "  +6.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.350.00"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.69"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.052.59"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.00"
$ws.Range("E17").Value = "  +14.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.201.34"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.80"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  -4.77%  "
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.89"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -5.09%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.46"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.93"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.01"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0617"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0893"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.24"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.19"
$ws.Range("E39").Value = "  +13.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.20"
$ws.Range("E40").Value = "  +15.02%  "
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.74"
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0933"
$ws.Range("E45").Value = "  -21.72%  "
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.272.67"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.80"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.233.09"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.45"
$ws.Range("E51").Value = "  -1.65%  "
